$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6669418
$ws.Range("J17").Value = 7145669
$ws.Range("L17").Value = 21437007
$ws.Range("N17").Value = -21437343

$ws.Range("H40").Value = 1263.5
$ws.Range("I40").Value = 750
$ws.Range("J40").Value = 1571.6
$ws.Range("K40").Value = 750
$ws.Range("L40").Value = 1571.6
$ws.Range("M40").Value = -575
$ws.Range("N40").Value = -1921.6

$ws.Range("H92").Value = 563.1429000000001
$ws.Range("I92").Value = 619.9
$ws.Range("J92").Value = 421.25
$ws.Range("K92").Value = 619.9
$ws.Range("L92").Value = 421.25
$ws.Range("M92").Value = 628.1
$ws.Range("N92").Value = -2917.25

$ws.Range("H113").Value = 333338020
$ws.Range("I113").Value = 500001000
$ws.Range("K113").Value = 500001000
$ws.Range("M113").Value = -499997746

$ws.Range("H129").Value = 147933.69
$ws.Range("J129").Value = 176435.48
$ws.Range("L129").Value = 529306.4400000001
$ws.Range("N129").Value = -539306.4400000001

$ws.Range("H132").Value = 2748.7837
$ws.Range("I132").Value = 3027.5938
$ws.Range("J132").Value = 964.4
$ws.Range("K132").Value = 9082.7814
$ws.Range("L132").Value = 2893.2
$ws.Range("M132").Value = -6552.7814
$ws.Range("N132").Value = -7953.2

$ws.Range("H137").Value = 1291.3462
$ws.Range("I137").Value = 1251.75
$ws.Range("J137").Value = 1423.3334
$ws.Range("K137").Value = 3755.25
$ws.Range("L137").Value = 4270.0002
$ws.Range("M137").Value = -1205.25
$ws.Range("N137").Value = -9370.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 669.8125
$ws.Range("I22").Value = 624.38464
$ws.Range("K22").Value = 624.38464
$ws.Range("M22").Value = -451.38464

$ws.Range("H105").Value = 1853616.6
$ws.Range("I105").Value = 1531.125
$ws.Range("K105").Value = 1531.125
$ws.Range("M105").Value = 215.875

$ws.Range("H134").Value = 2806.0208
$ws.Range("I134").Value = 2842.319
$ws.Range("K134").Value = 8526.957
$ws.Range("M134").Value = -5991.957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.27274
$ws.Range("I7").Value = 19
$ws.Range("K7").Value = 19
$ws.Range("M7").Value = 94

$ws.Range("H31").Value = 2867.6667
$ws.Range("I31").Value = 1435.6111
$ws.Range("J31").Value = 6304.6
$ws.Range("K31").Value = 1435.6111
$ws.Range("L31").Value = 6304.6
$ws.Range("M31").Value = -1140.6111
$ws.Range("N31").Value = -6894.6

$ws.Range("H34").Value = 2867.6667
$ws.Range("I34").Value = 1435.6111
$ws.Range("J34").Value = 6304.6
$ws.Range("K34").Value = 1435.6111
$ws.Range("L34").Value = 6304.6
$ws.Range("M34").Value = -1233.6111
$ws.Range("N34").Value = -6708.6

$ws.Range("H105").Value = 1037.4
$ws.Range("I105").Value = 807.875
$ws.Range("K105").Value = 807.875
$ws.Range("M105").Value = 939.125

$ws.Range("H110").Value = 30351
$ws.Range("J110").Value = 30351
$ws.Range("L110").Value = 30351
$ws.Range("N110").Value = -38531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1248.6666
$ws.Range("I5").Value = 1004.32
$ws.Range("J5").Value = 1804
$ws.Range("K5").Value = 3012.96
$ws.Range("L5").Value = 5412
$ws.Range("M5").Value = -2900.96
$ws.Range("N5").Value = -5636

$ws.Range("H7").Value = 230.8
$ws.Range("I7").Value = 99
$ws.Range("K7").Value = 297
$ws.Range("M7").Value = -185

$ws.Range("H34").Value = 618.9231
$ws.Range("J34").Value = 844.1111
$ws.Range("L34").Value = 2532.3333
$ws.Range("N34").Value = -2700.3333

$ws.Range("H55").Value = 2493
$ws.Range("I55").Value = 90
$ws.Range("J55").Value = 2862.6924
$ws.Range("K55").Value = 270
$ws.Range("L55").Value = 8588.0772
$ws.Range("M55").Value = -93
$ws.Range("N55").Value = -8942.0772

$ws.Range("H75").Value = 1717.091
$ws.Range("I75").Value = 971
$ws.Range("J75").Value = 1996.875
$ws.Range("K75").Value = 2913
$ws.Range("L75").Value = 5990.625
$ws.Range("M75").Value = -1915
$ws.Range("N75").Value = -7986.625

$ws.Range("H78").Value = 1717.091
$ws.Range("I78").Value = 971
$ws.Range("J78").Value = 1996.875
$ws.Range("K78").Value = 8739
$ws.Range("L78").Value = 17971.875
$ws.Range("M78").Value = -3747
$ws.Range("N78").Value = -27955.875

$ws.Range("H80").Value = 27000.25
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 27000.25
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 81000.75
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -82872.75

$ws.Range("H83").Value = 27000.25
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 27000.25
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 243002.25
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -252362.25

$ws.Range("H100").Value = 5639.8184
$ws.Range("J100").Value = 5639.8184
$ws.Range("L100").Value = 16919.4552
$ws.Range("N100").Value = -18541.4552

$ws.Range("H107").Value = 6100.9414
$ws.Range("I107").Value = 7344.2856
$ws.Range("J107").Value = 298.66666
$ws.Range("K107").Value = 22032.8568
$ws.Range("L107").Value = 895.9999799999999
$ws.Range("M107").Value = -20112.8568
$ws.Range("N107").Value = -4735.99998

$ws.Range("H117").Value = 1171.9286
$ws.Range("J117").Value = 1088.909
$ws.Range("L117").Value = 3266.727
$ws.Range("N117").Value = -10150.727

$ws.Range("H122").Value = 674.2381
$ws.Range("I122").Value = 448
$ws.Range("J122").Value = 764.73334
$ws.Range("K122").Value = 4032
$ws.Range("L122").Value = 6882.60006
$ws.Range("M122").Value = -1582
$ws.Range("N122").Value = -11782.60006

$ws.Range("H131").Value = 675.385
$ws.Range("J131").Value = 692.5645
$ws.Range("L131").Value = 2077.6935
$ws.Range("N131").Value = -12157.6935

$ws.Range("H135").Value = 1248.6666
$ws.Range("I135").Value = 1004.32
$ws.Range("J135").Value = 1804
$ws.Range("K135").Value = 9038.880000000001
$ws.Range("L135").Value = 16236
$ws.Range("M135").Value = -6503.880000000001
$ws.Range("N135").Value = -21306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1543.6957
$ws.Range("I97").Value = 1257.579
$ws.Range("J97").Value = 2902.75
$ws.Range("K97").Value = 1257.579
$ws.Range("L97").Value = 2902.75
$ws.Range("M97").Value = -761.579
$ws.Range("N97").Value = -3894.75

$ws.Range("H102").Value = 1669.1464
$ws.Range("I102").Value = 1448.8788
$ws.Range("K102").Value = 1448.8788
$ws.Range("M102").Value = 173.1212

$ws.Range("H113").Value = 14707.286
$ws.Range("I113").Value = 19000.2
$ws.Range("J113").Value = 3975
$ws.Range("K113").Value = 19000.2
$ws.Range("L113").Value = 3975
$ws.Range("M113").Value = -16830.2
$ws.Range("N113").Value = -8315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2427.5945
$ws.Range("I40").Value = 2331.2
$ws.Range("K40").Value = 2331.2
$ws.Range("M40").Value = -2195.2

$ws.Range("H61").Value = 5509.25
$ws.Range("I61").Value = 1790.2222
$ws.Range("J61").Value = 16666.334
$ws.Range("K61").Value = 1790.2222
$ws.Range("L61").Value = 16666.334
$ws.Range("M61").Value = -1588.2222
$ws.Range("N61").Value = -17070.334

$ws.Range("H113").Value = 5509.25
$ws.Range("I113").Value = 1790.2222
$ws.Range("J113").Value = 16666.334
$ws.Range("K113").Value = 1790.2222
$ws.Range("L113").Value = 16666.334
$ws.Range("M113").Value = 379.7778000000001
$ws.Range("N113").Value = -21006.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3103.7144
$ws.Range("J81").Value = 3545.2
$ws.Range("L81").Value = 7090.4
$ws.Range("N81").Value = -9212.4

$ws.Range("H84").Value = 3103.7144
$ws.Range("J84").Value = 3545.2
$ws.Range("L84").Value = 35452
$ws.Range("N84").Value = -46060

$ws.Range("H113").Value = 1078.875
$ws.Range("I113").Value = 1296.9231
$ws.Range("J113").Value = 134
$ws.Range("K113").Value = 3890.7693
$ws.Range("L113").Value = 402
$ws.Range("M113").Value = -1720.7693
$ws.Range("N113").Value = -4742

$ws.Range("H132").Value = 1254.7179
$ws.Range("I132").Value = 794.1429000000001
$ws.Range("J132").Value = 2427.0908
$ws.Range("K132").Value = 2382.4287
$ws.Range("L132").Value = 7281.2724
$ws.Range("M132").Value = 147.5712999999996
$ws.Range("N132").Value = -12341.2724
